$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.457.41"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "2.931.25"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'594.73"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").Value = "'145.09"
$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'0.501"
$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("D9").Value = "'6.99"
$ws.Range("E9").Value = "  +1.22%  "

$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = "  -1.75%  "

$ws.Range("D11").Value = "'0.439"
$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").Value = "'0.0000224"
$ws.Range("E12").Value = "  -0.86%  "

$ws.Range("D13").Value = "'33.22"
$ws.Range("E13").Value = "  -1.19%  "

$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").Value = "3.417.68"
$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("D16").Value = "61.438.68"
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("D17").Value = "2.935.80"
$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("D18").Value = "'6.65"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").Value = "'434.15"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").Value = "'13.58"
$ws.Range("E20").Value = "  +1.49%  "

$ws.Range("D21").Value = "'0.673"
$ws.Range("E21").Value = "  -0.87%  "

$ws.Range("D22").Value = "'7.08"
$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("D23").Value = "'81.64"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").Value = "'10.99"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("E25").Value = "  -1.13%  "

$ws.Range("D26").Value = "'11.72"
$ws.Range("E26").Value = "  -0.96%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("E28").Value = "  -4.72%  "

$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("D30").Value = "'6.93"
$ws.Range("E30").Value = "  -1.37%  "

$ws.Range("D31").Value = "'26.68"
$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("D32").Value = "'0.108"
$ws.Range("E32").Value = "  +1.23%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").Value = "0.0₃0879"
$ws.Range("E34").Value = "  +2.02%  "

$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("E36").Value = "  -0.42%  "

$ws.Range("D37").Value = "'2.99"
$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("E39").Value = "  -0.40%  "

$ws.Range("D40").Value = "'8.47"
$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("D41").Value = "'42.11"
$ws.Range("E41").Value = "  +2.45%  "

$ws.Range("D42").Value = "'0.278"
$ws.Range("E42").Value = "  -3.37%  "

$ws.Range("D43").Value = "'0.0343"
$ws.Range("E43").Value = "  -0.36%  "

$ws.Range("D44").Value = "2.687.18"
$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'133.85"
$ws.Range("E45").Value = "  +0.40%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'362.24"
$ws.Range("E46").Value = "  -3.79%  "

$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("D48").Value = "'23.46"
$ws.Range("E48").Value = "  -1.47%  "

$ws.Range("E49").Value = "  -0.91%  "

$ws.Range("D50").Value = "'2.01"

$ws.Range("E51").Value = "  +0.50%  "
